$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $range as a genuine text (shared-string) cell,
# without Excel's "looks like a number" auto-conversion turning pure-digit
# strings (e.g. "20093135", "1", "2" ...) into numeric cells, and without
# touching the cell's existing style/format. We do this by putting a
# literal-string formula ("=""...""") in a scratch cell, copying it, and
# pasting *values only* into the destination - the pasted result keeps the
# destination's original formatting but is stored as text.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""$escaped"""
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# Row 12 is new (the sheet previously only went down to row 11). Give it the
# same look (borders/style) as the row above before filling it in, by
# copying the formatting only.
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Replace the product listing (rows 2-11) and add the new row (12) with the
# updated data.
Set-TextValue $ws.Range("A2")  '10003485'
Set-TextValue $ws.Range("B2")  'S/Q CHOCO ALMOND 52G'
Set-TextValue $ws.Range("C2")  'PM1MKT'
Set-TextValue $ws.Range("D2")  '1'
Set-TextValue $ws.Range("E2")  '1'
Set-TextValue $ws.Range("F2")  'RT,(E-1B)'

Set-TextValue $ws.Range("A3")  '10036987'
Set-TextValue $ws.Range("B3")  'S/Q CHOCO CASHEW 52G'
Set-TextValue $ws.Range("C3")  'PM1MKT'
Set-TextValue $ws.Range("D3")  '1'
Set-TextValue $ws.Range("E3")  '2'
Set-TextValue $ws.Range("F3")  'RT,(E-1B)'

Set-TextValue $ws.Range("A4")  '10003814'
Set-TextValue $ws.Range("B4")  'AQUA AIR MINERAL 600'
Set-TextValue $ws.Range("C4")  'PM1MKT'
Set-TextValue $ws.Range("D4")  '2'
Set-TextValue $ws.Range("E4")  '1'
Set-TextValue $ws.Range("F4")  'RT,(E-3B)'

Set-TextValue $ws.Range("A5")  '20040383'
Set-TextValue $ws.Range("B5")  'NU MILK TEA 330ML'
Set-TextValue $ws.Range("C5")  'PM1MKT'
Set-TextValue $ws.Range("D5")  '2'
Set-TextValue $ws.Range("E5")  '2'
Set-TextValue $ws.Range("F5")  'RT,(E-1B)'

Set-TextValue $ws.Range("A6")  '20069527'
Set-TextValue $ws.Range("B6")  'NU TEH TARIK 330ML'
Set-TextValue $ws.Range("C6")  'PM1MKT'
Set-TextValue $ws.Range("D6")  '2'
Set-TextValue $ws.Range("E6")  '3'
Set-TextValue $ws.Range("F6")  'RT,(E-1B)'

Set-TextValue $ws.Range("A7")  '10014404'
Set-TextValue $ws.Range("B7")  'MILO HEALTY DRINK220'
Set-TextValue $ws.Range("C7")  'PM1MKT'
Set-TextValue $ws.Range("D7")  '2'
Set-TextValue $ws.Range("E7")  '4'
Set-TextValue $ws.Range("F7")  'RT,(E-2B)'

Set-TextValue $ws.Range("A8")  '20068536'
Set-TextValue $ws.Range("B8")  'REBO KUACI G.TEA 120'
Set-TextValue $ws.Range("C8")  'PM1MKT'
Set-TextValue $ws.Range("D8")  '3'
Set-TextValue $ws.Range("E8")  '1'
Set-TextValue $ws.Range("F8")  'RT,(E-1B)'

Set-TextValue $ws.Range("A9")  '20098334'
Set-TextValue $ws.Range("B9")  'REBO KUACI CRM 120G'
Set-TextValue $ws.Range("C9")  'PM1MKT'
Set-TextValue $ws.Range("D9")  '3'
Set-TextValue $ws.Range("E9")  '2'
Set-TextValue $ws.Range("F9")  'RT,(E-1B)'

Set-TextValue $ws.Range("A10") '20092331'
Set-TextValue $ws.Range("B10") 'EKONOMI LIQ JR.NP650'
Set-TextValue $ws.Range("C10") 'PM1MKT'
Set-TextValue $ws.Range("D10") '4'
Set-TextValue $ws.Range("E10") '1'
Set-TextValue $ws.Range("F10") 'RT,(E-1B)'

Set-TextValue $ws.Range("A11") '20011008'
Set-TextValue $ws.Range("B11") 'LIFEBUOY BW RED 400'
Set-TextValue $ws.Range("C11") 'PM1MKT'
Set-TextValue $ws.Range("D11") '4'
Set-TextValue $ws.Range("E11") '2'
Set-TextValue $ws.Range("F11") 'PT,(E-3B)'

Set-TextValue $ws.Range("A12") '20040194'
Set-TextValue $ws.Range("B12") 'GRNR BC FWS.VT.C100'
Set-TextValue $ws.Range("C12") 'PM1MKT'
Set-TextValue $ws.Range("D12") '4'
Set-TextValue $ws.Range("E12") '3'
Set-TextValue $ws.Range("F12") 'RT,(E-2B)'

# Clean up the scratch cell used for the text-forcing trick.
$ws.Range("ZZ1").ClearContents() | Out-Null
